# gsc-export/HTTPS.xlsx -- roll the reporting window forward by two days.
#
# The "Chart" sheet lists one row per date (column A) together with the
# Non-HTTPS (B) and HTTPS (C) page counts for that date. The export was
# refreshed: the two oldest days (2025-11-11 and 2025-11-12) drop off the
# front of the window, every remaining row shifts up by two, and two new
# days (2026-02-10, 2026-02-11) are appended at the end with updated
# HTTPS page counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Chart" sheet

# Drop the two oldest date rows (2025-11-11, 2025-11-12). Excel shifts all
# rows below them up by two, which reproduces the "every value moves up by
# two rows" pattern seen throughout the diff.
$ws.Rows("2:3").Delete()

# After the shift the last populated row is 90 (2026-02-09). Append the two
# new trailing days. A literal date-formatted string assigned straight into
# Value2 would be auto-converted into a date serial by Excel's type
# inference, so instead we stage the text in a scratch cell via a text
# formula (never auto-converted) and copy/paste just the value across -
# this keeps column A stored as plain text, matching the rest of the sheet.
$helper = $ws.Range("Z1")

$helper.Formula = "=""2026-02-10"""
$helper.Copy()
$ws.Range("A91").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$helper.Clear()

$helper.Formula = "=""2026-02-11"""
$helper.Copy()
$ws.Range("A92").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$helper.Clear()

$ws.Range("B91").Value2 = 0.0
$ws.Range("C91").Value2 = 29.0

$ws.Range("B92").Value2 = 0.0
$ws.Range("C92").Value2 = 29.0
